# Summer 23 week 13 inputs - append new matchup rows to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 1071

$data = @(
    @(5,14,6,6),
    @(4,8,3,12),
    @(4,13,5,7),
    @(4,7,5,13),
    @(3,15,4,5),
    @(4,5,2,15),
    @(5,5,8,15),
    @(4,17,3,3),
    @(4,15,5,5),
    @(3,7,4,13),
    @(3,2,5,18),
    @(9,17,7,3),
    @(3,13,4,7),
    @(5,15,3,5),
    @(6,7,5,13),
    @(6,16,4,4),
    @(3,13,4,7),
    @(3,12,4,8),
    @(6,14,4,6),
    @(4,17,3,3),
    @(2,3,4,17),
    @(5,8,4,12),
    @(4,14,2,6),
    @(5,15,3,5),
    @(5,16,4,4),
    @(9,12,4,8),
    @(3,14,6,6),
    @(3,13,4,7)
)

$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view to match the committed state: scroll/select the next empty row
$newActiveRow = $endRow + 1
$ws.Application.ActiveWindow.ScrollRow = $startRow + 11
$ws.Range("A" + $newActiveRow).Select()
